# Insert a new data row at row 150 (pushing the existing rows 150-236
# down to 151-237) and populate it with the new "Frutilla" price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(150).Insert()

$ws.Cells.Item(150, 1).Value  = 10
$ws.Cells.Item(150, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(150, 3).Value  = "La Araucanía"
$ws.Cells.Item(150, 4).Value  = 44582
$ws.Cells.Item(150, 5).Value  = 9
$ws.Cells.Item(150, 6).Value  = "Fruta"
$ws.Cells.Item(150, 7).Value  = 100101
$ws.Cells.Item(150, 8).Value  = "Berries"
$ws.Cells.Item(150, 9).Value  = 100112025
$ws.Cells.Item(150, 10).Value = "Frutilla"
$ws.Cells.Item(150, 11).Value = "Sin especificar"
$ws.Cells.Item(150, 12).Value = "Primera"
$ws.Cells.Item(150, 13).Value = 110
$ws.Cells.Item(150, 14).Value = 7000
$ws.Cells.Item(150, 15).Value = 7000
$ws.Cells.Item(150, 16).Value = 7000
$ws.Cells.Item(150, 17).Value = "`$/caja 7 kilos"
$ws.Cells.Item(150, 18).Value = "Región de La Araucanía"
$ws.Cells.Item(150, 19).Value = 1000
$ws.Cells.Item(150, 20).Value = 7
